# Applies the cryptos list refresh described in the commit
# "Updated cryptos list on Fri Jun 23 06:52:46 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values look like numbers (e.g. "0.9998", "0.06640",
# "29.989.56"). Excel auto-parses a plain string assignment that looks
# numeric and would silently turn it into a real number (destroying
# trailing zeros / thousand-dot formatting). Force the cell to Text
# format first so the literal string is preserved, matching the source
# data which stores these as inline strings.
$dUpdates = @{
    "D2" = "29.989.56"
    "D3" = "1.879.46"
    "D4" = "0.9998"
    "D5" = "243.00"
    "D6" = "0.9996"
    "D7" = "0.4955"
    "D8" = "0.2921"
    "D9" = "0.06640"
    "D10" = "1.880.53"
    "D11" = "16.78"
    "D12" = "0.07250"
    "D13" = "0.6672"
    "D14" = "86.55"
    "D15" = "4.916"
    "D16" = "29.967.01"
    "D17" = "0.000007874"
    "D18" = "0.9991"
    "D19" = "12.80"
    "D20" = "2.122.75"
    "D21" = "0.9989"
    "D22" = "4.778"
    "D23" = "5.749"
    "D24" = "9.076"
    "D25" = "142.50"
    "D26" = "149.60"
    "D29" = "1.393"
    "D30" = "4.199"
    "D31" = "0.08754"
    "D32" = "3.970"
    "D33" = "0.05072"
    "D34" = "1.117"
    "D35" = "0.7121"
    "D36" = "2.669"
    "D38" = "2.691"
    "D39" = "2.177"
    "D40" = "0.9337"
    "D41" = "5.802"
    "D42" = "0.4249"
    "D43" = "0.9988"
    "D44" = "102.62"
    "D45" = "7.451"
    "D46" = "0.1269"
    "D47" = "0.05668"
    "D48" = "32.61"
    "D49" = "8.327"
    "D50" = "0.3788"
    "D51" = "55.99"
}
foreach ($addr in $dUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dUpdates[$addr]
}

# Coin name / link / volume columns are plain text and safe to set directly.
$otherUpdates = @{
    "E2" = "  -0.40%  "
    "E3" = "  -1.58%  "
    "E4" = "  -0.11%  "
    "E5" = "  -3.30%  "
    "E6" = "  -0.14%  "
    "E7" = "  -2.66%  "
    "E8" = "  -1.08%  "
    "E10" = "  -1.63%  "
    "E11" = "  -2.63%  "
    "E12" = "  -1.51%  "
    "E13" = "  -3.12%  "
    "E14" = "  +0.17%  "
    "E15" = "  +0.94%  "
    "E16" = "  -0.49%  "
    "E17" = "  -2.75%  "
    "E18" = "  -0.17%  "
    "E19" = "  -1.29%  "
    "E20" = "  -1.72%  "
    "E21" = "  -0.18%  "
    "E22" = "  -0.94%  "
    "E23" = "  +0.24%  "
    "E24" = "  -0.69%  "
    "E25" = "  +5.35%  "
    "E26" = "  +1.86%  "
    "E27" = "  +0.06%  "
    "E28" = "  -3.66%  "
    "E29" = "  +0.16%  "
    "E30" = "  -0.57%  "
    "E31" = "  -0.33%  "
    "E32" = "  -0.88%  "
    "E33" = "  +0.23%  "
    "E34" = "  -2.21%  "
    "E35" = "  -0.05%  "
    "E36" = "  -0.85%  "
    "E37" = "  +6.72%  "
    "E38" = "  -4.15%  "
    "E39" = "  -4.03%  "
    "E40" = "  -3.68%  "
    "E41" = "  -5.30%  "
    "E42" = "  -0.77%  "
    "E43" = "  -0.02%  "
    "E44" = "  -1.97%  "
    "E45" = "  -1.75%  "
    "E46" = "  -0.55%  "
    "E47" = "  -1.20%  "
    "E48" = "  -1.33%  "
    "B49" = "EnergySwap"
    "C49" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "E49" = "  -1.34%  "
    "B50" = "Decentraland"
    "C50" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
    "E50" = "  -0.11%  "
    "E51" = "  -1.38%  "
}
foreach ($addr in $otherUpdates.Keys) {
    $ws.Range($addr).Value = $otherUpdates[$addr]
}
